$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F ("想去人数") that changed.
$updates = @{
    2  = 1132
    8  = 2364
    9  = 7706
    10 = 918
    11 = 438
    12 = 374
    14 = 417
    16 = 7887
    18 = 1372
    24 = 156
    35 = 44
}

# Both "展览" and "全部类型" sheets contain the same rows of data and need
# the same F-column updates applied.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
